$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" caption in A1
$ws.Range("A1").Value = "Datos actualizados a 28 de Julio de 2020 a las 07:35"

# Row 6: India -> India
$ws.Range("B6").Value = 1484136
$ws.Range("C6").Value = 1633
$ws.Range("D6").Value = 954004
$ws.Range("E6").Value = 496671
$ws.Range("G6").Value = 13
$ws.Range("H6").Value = 33461

# Row 15: Pakistan -> Pakistan
$ws.Range("B15").Value = 275225
$ws.Range("C15").Value = 936
$ws.Range("D15").Value = 242436
$ws.Range("E15").Value = 26924
$ws.Range("G15").Value = 23
$ws.Range("H15").Value = 5865

# Row 22: Francia -> Francia
$ws.Range("D22").Value = 81082
$ws.Range("E22").Value = 71788

# Row 29: China -> Kazajistan
$ws.Range("A29").Value = "Kazajistan"
$ws.Range("B29").Value = 84648
$ws.Range("C29").Value = 1526
$ws.Range("D29").Value = 54404
$ws.Range("E29").Value = 29659
$ws.Range("H29").Value = 585

# Row 30: Kazajistan -> China
$ws.Range("A30").Value = "China"
$ws.Range("B30").Value = 83959
$ws.Range("C30").Value = 68
$ws.Range("D30").Value = 78934
$ws.Range("E30").Value = 391
$ws.Range("H30").Value = 4634

# Row 56: Ghana -> Kirguistan
$ws.Range("A56").Value = "Kirguistan"
$ws.Range("B56").Value = 33718
$ws.Range("C56").Value = 422
$ws.Range("D56").Value = 22296
$ws.Range("E56").Value = 10093
$ws.Range("G56").Value = 28
$ws.Range("H56").Value = 1329

# Row 57: Kirguistan -> Ghana
$ws.Range("A57").Value = "Ghana"
$ws.Range("B57").Value = 33624
$ws.Range("D57").Value = 29801
$ws.Range("E57").Value = 3655
$ws.Range("H57").Value = 168

# Row 64: Uzbekistan -> Uzbekistan
$ws.Range("B64").Value = 21506
$ws.Range("C64").Value = 297
$ws.Range("E64").Value = 9710
$ws.Range("G64").Value = 1
$ws.Range("H64").Value = 122

# Row 109: Tailandia -> Tailandia
$ws.Range("B109").Value = 3297
$ws.Range("C109").Value = 2
$ws.Range("E109").Value = 128

# Row 144: Niger -> Niger
$ws.Range("B144").Value = 1132
$ws.Range("E144").Value = 36

# Row 189: Butan -> Butan
$ws.Range("B189").Value = 99
$ws.Range("C189").Value = 4
$ws.Range("D189").Value = 86
$ws.Range("E189").Value = 13
